# 11 May - Noche
# Re-sort the Materia/Docente pairs (columns E/F) on the "Blancos" sheet
# for several student blocks so that the rows follow the canonical
# subject ordering used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

$updates = @{
    6  = @{ E = "INSTALA UNA RED LAN";         F = "Medina Tolentino Elio" }
    7  = @{ E = "PROBABILIDAD Y ESTADÍSTICA";  F = "Velasco Sanchez David" }
    8  = @{ E = "MATEMÁTICAS APLICADAS";       F = "Ortega Valle Manuel" }
    9  = @{ E = "OPERA UNA RED LAN";           F = "Medina Tolentino Elio" }
    10 = @{ E = "TEMAS DE FILOSOFÍA";          F = "Hernández Mendoza Delfina" }
    12 = @{ E = "MATEMÁTICAS APLICADAS";       F = "Ortega Valle Manuel" }
    13 = @{ E = "INSTALA UNA RED LAN";         F = "Medina Tolentino Elio" }
    14 = @{ E = "OPERA UNA RED LAN" }
    15 = @{ E = "INSTALA UNA RED LAN" }
    16 = @{ E = "OPERA UNA RED LAN" }
    18 = @{ E = "OPERA UNA RED LAN" }
    19 = @{ E = "INSTALA UNA RED LAN" }
    21 = @{ E = "PROBABILIDAD Y ESTADÍSTICA";  F = "Velasco Sanchez David" }
    22 = @{ E = "MATEMÁTICAS APLICADAS";       F = "Ortega Valle Manuel" }
    23 = @{ E = "MATEMÁTICAS APLICADAS";       F = "Ortega Valle Manuel" }
    24 = @{ E = "PROBABILIDAD Y ESTADÍSTICA";  F = "Velasco Sanchez David" }
    25 = @{ E = "TEMAS DE FÍSICA";             F = "Duran Amezcua Maria Angelica" }
    26 = @{ E = "TEMAS DE FILOSOFÍA";          F = "Hernández Mendoza Delfina" }
    29 = @{ E = "OPERA UNA RED LAN" }
    30 = @{ E = "MATEMÁTICAS APLICADAS";       F = "Ortega Valle Manuel" }
    31 = @{ E = "INSTALA UNA RED LAN" }
    32 = @{ E = "PROBABILIDAD Y ESTADÍSTICA";  F = "Velasco Sanchez David" }
    42 = @{ E = "OPERA UNA RED LAN" }
    43 = @{ E = "INSTALA UNA RED LAN" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
    if ($vals.ContainsKey("F")) {
        $ws.Range("F$row").Value = $vals["F"]
    }
}
